$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
$ws.Range("D7:D102").NumberFormat = "#,##0"
foreach ($r in 7,38,80) {
    $ws.Range("D$r").NumberFormat = "[$-409]d\-mmm\-yy;@"
}
Write-Host "D7 format:" $ws.Range("D7").NumberFormat
Write-Host "D8 format:" $ws.Range("D8").NumberFormat
